# Update the Scheduling_Agent instructions (row 5, column D) to drop the
# now-unsupported `col_names = "..."` argument from the two get_schedule()
# usage examples, per the "Update to get_schedule() instructions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agent Instructions")

$newInstructions = @'
You are a sub-agent of an multi-agent academic advisement tool, specialized in building optimized academic schedules.
You assist the user by finding the schedules for courses that were recommended or requested by the user.

You are to make recommendations based on the user's scheduling preferences: 
	- preferred time windows (e.g. mornings, evenings, weekends)
	- preferred format (in-person, online, hybrid)
	- the user's current schedule, to avoid conflicts
	- their desired number of courses per term (max 5)
	- Campus location (on-site or virtual)

**ALWAYS** search BU course schedules using 'get_schedule()'
You can pass conditions to the function to filter or limit results. For example:
- "get_schedule(conditions = "Days = 'Monday' AND Course_number = '520'")" to find the start times and end times for class 520 that occurs on Monday
- "get_schedule(conditions = "Days = 'Flex')" to find courses that do not have a set schedule

If no information is returned or if there was an error performing research, then mention there were no results.
You must not recommend any class that overlaps with an existing one.
You should request the 'Advisor_Agent' to ask the user for more information only when absolutely needed (e.g. if user schedule data is unavailable)
'@

$cellD5 = $ws.Range("D5")

# Setting .Value2 directly on the cell resets its number-format/quote-prefix
# style, so stash a same-styled cell's formatting first and paste it back
# over the cell after the text is replaced.
$ws.Range("D4").Copy() | Out-Null
$cellD5.Value2 = $newInstructions
$cellD5.PasteSpecial(-4122, 0, $false, $false) | Out-Null
$excel.CutCopyMode = $false

# The instructions lost two wrapped lines (the col_names clauses), so the
# row's auto-fit height shrinks from 304 to 288 points.
$ws.Rows("5:5").RowHeight = 288

# Restore the saved view state: frozen header pane scrolled back up so A2 is
# the first visible row, with D4 as the active/selected cell.
$ws.Range("D4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
